$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the existing hyperlinks and wipe the current data rows (2-7),
# keeping the header row (row 1) untouched.
$ws.Hyperlinks.Delete()
$ws.Range("A2:I7").Clear()

# Row 2: Doru3 / Student
$ws.Range("A2").Value = "bocaioandoru12+3@gmail.com"
$ws.Range("B2").Value = "Doru3"
$ws.Range("C2").Value = "Student"
$ws.Range("D2").Value = "UTCN"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = "CTI"

# Row 3: Doru4 / Campus_Student
$ws.Range("A3").Value = "bocaioandoru12+4@gmail.com"
$ws.Range("B3").Value = "Doru4"
$ws.Range("C3").Value = "Campus_Student"
$ws.Range("D3").Value = "UTCN"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = "CTI"

# Row 4: Doru2 / Professor
$ws.Range("A4").Value = "bocaioandoru12+2@gmail.com"
$ws.Range("B4").Value = "Doru2"
$ws.Range("C4").Value = "Professor"
$ws.Range("D4").Value = "UTCN"
$ws.Range("H4").Value = "CTI"
$ws.Range("I4").Value = "Eng."

# Row 6: keeps the hyperlink styling of the old table but carries no value.
$ws.Range("A6").Style = "Hyperlink"

# Row 7: restore the originally empty, hyperlink-styled placeholder cell.
$ws.Range("A7").Style = "Hyperlink"

# Hyperlinks for the new table (targets keep the same rId mapping order as
# before the edit, so the mailto addresses are intentionally swapped versus
# the cell text, exactly like the source workbook).
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:bocaioandoru12+4@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:bocaioandoru12+3@gmail.com")

# Re-apply the Hyperlink cell style so A2/A3 keep the same style index the
# workbook used before the edit.
$ws.Range("A2").Style = "Hyperlink"
$ws.Range("A3").Style = "Hyperlink"

# Update the active selection to match the saved view state.
$ws.Range("H12").Select()
